$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(22, 3)
$cell.Value = "'318"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.Value = "'923503.99"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 3)
$cell.Value = "'110"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'423337.00"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 3)
$cell.Value = "'36"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'170233.00"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 3)
$cell.Value = "'86"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 4)
$cell.Value = "'225156.00"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 3)
$cell.Value = "'488"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'1432983.41"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 3)
$cell.Value = "'195"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'907647.11"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 3)
$cell.Value = "'23"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'140500.00"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 3)
$cell.Value = "'88"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'242933.17"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 3)
$cell.Value = "'522"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'1660527.52"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(53, 3)
$cell.Value = "'78"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(53, 4)
$cell.Value = "'432878.23"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(80, 3)
$cell.Value = "'830"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(80, 4)
$cell.Value = "'2502044.56"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(81, 3)
$cell.Value = "'308"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(81, 4)
$cell.Value = "'1166766.79"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(83, 3)
$cell.Value = "'26"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(83, 4)
$cell.Value = "'149080.04"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(91, 3)
$cell.Value = "'91"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(91, 4)
$cell.Value = "'219878.00"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(92, 3)
$cell.Value = "'387"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(92, 4)
$cell.Value = "'1108260.67"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(93, 3)
$cell.Value = "'165"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(93, 4)
$cell.Value = "'640292.91"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(94, 3)
$cell.Value = "'46"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(94, 4)
$cell.Value = "'211347.01"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(96, 3)
$cell.Value = "'11"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(96, 4)
$cell.Value = "'22000.00"
$cell.Style = "Normal"

